# Update countries & provincias Spain
# Applies the COVID-19 "Pais" dashboard refresh: updated case counters for a
# handful of countries, Guyana overtaking Botsuana in the ranking (rows 146/147
# swap countries while keeping their row position), and a refreshed timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Header timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 23:16"

# --- Updated country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
Set-Row 4 8029238 37240 5166131 2643135 0 277 219972

# Row 6: Brasil
Set-Row 6 5103408 8429 4470165 482554 0 183 150689

# Row 25: Alemania
Set-Row 25 331093 4802 276900 44472 0 19 9721

# Row 29: Canada
Set-Row 29 182791 927 154246 18918 0 14 9627

# Row 103: Namibia
Set-Row 103 11989 53 9913 1947 0 1 129

# Row 105: Maldivas
Set-Row 105 10943 49 9783 1125 0 0 35

# Row 118: Cabo Verde
Set-Row 118 7155 83 6075 1005 0 0 75

# Row 136: Siria
Set-Row 136 4774 56 1331 3215 0 4 228

# Row 144: Gambia
Set-Row 144 3636 4 2593 925 0 1 118

# Rows 146/147: Guyana overtakes Botsuana in the ranking, so the country names
# in column A swap while the row position (and thus the rank number in A)
# stays put; Guyana gets freshly updated figures, Botsuana keeps its previous
# (unchanged) figures one row further down.
$ws.Range("A146").Value = "Guyana"
Set-Row 146 3521 52 2391 1026 0 1 104

$ws.Range("A147").Value = "Botsuana"
Set-Row 147 3515 296 853 2642 0 2 20

# Row 174: Curazao
Set-Row 174 585 2 332 252 0 0 1
